$wb = $excel.ActiveWorkbook

# Use the existing "UTI" sheet's header formatting as the style source
$uti = $wb.Worksheets.Item("UTI")
$srcHeader = $uti.Range("A1:D1")

# Add the new "UPA" sheet after the last existing sheet ("UTI")
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "UPA"

# Header row
$ws.Range("A1").Value = "LEITO"
$ws.Range("B1").Value = "NOME DO PACIENTE"
$ws.Range("C1").Value = "DIETA"
$ws.Range("D1").Value = "OBSERVAÇÕES"

# Match the bold / bordered / centered header style used on the other sheets
$srcHeader.Copy()
$ws.Range("A1:D1").PasteSpecial(-4122)

# Data rows
$ws.Range("A2").Value = "ESTAB 01"
$ws.Range("B2").Value = "MICHAEL"
$ws.Range("C2").Value = "LIQUIDA"

$ws.Range("A3").Value = "ESTAB 02"
$ws.Range("B3").Value = "TREVOR"
$ws.Range("C3").Value = "PASTOSA"

$ws.Range("A4").Value = "ESTAB 03"
$ws.Range("B4").Value = "FRANKLIN"
$ws.Range("C4").Value = "LIVRE"

# Restore the originally active sheet/selection
[void]$wb.Worksheets.Item("Enfermaria").Activate()
[void]$wb.Worksheets.Item("Enfermaria").Range("A1").Select()
